$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 74

# Copy the formatting (style) of the previous data row onto the new row,
# so the new row's cells share the same style index as the rest of the log.
$ws.Range("G73").Copy()
$ws.Range("A74:H74").PasteSpecial(-4122)

$ws.Cells.Item($row, 1).Value = "2025-08-30 03:44:17 UTC"
$ws.Cells.Item($row, 2).Value = "2025-08-30 09:14:17 IST"
$ws.Cells.Item($row, 3).Value = "SKIPPED"
$ws.Cells.Item($row, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($row, 5).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-28-08-2025.pdf"
$ws.Cells.Item($row, 7).Value = 0
